$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 833847.2
$ws.Range("J6").Value = 1517
$ws.Range("L6").Value = 4551
$ws.Range("N6").Value = -4775
$ws.Range("H15").Value = 365.70587
$ws.Range("I15").Value = 365.70587
$ws.Range("K15").Value = 1097.11761
$ws.Range("M15").Value = -928.11761
$ws.Range("H33").Value = 332.5
$ws.Range("I33").Value = 332.5
$ws.Range("K33").Value = 332.5
$ws.Range("M33").Value = -103.5
$ws.Range("H41").Value = 370.26666
$ws.Range("I41").Value = 379.58334
$ws.Range("J41").Value = 333
$ws.Range("K41").Value = 379.58334
$ws.Range("L41").Value = 333
$ws.Range("M41").Value = 60.41665999999998
$ws.Range("N41").Value = -1213
$ws.Range("H115").Value = 1484.8
$ws.Range("I115").Value = 856
$ws.Range("K115").Value = 2568
$ws.Range("M115").Value = -1001
$ws.Range("H138").Value = 7114.3076
$ws.Range("I138").Value = 5615.1665
$ws.Range("J138").Value = 8399.286
$ws.Range("K138").Value = 16845.4995
$ws.Range("L138").Value = 25197.858
$ws.Range("M138").Value = -11705.4995
$ws.Range("N138").Value = -35477.858

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2666.6667
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 3500
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 3500
$ws.Range("M11").Value = -856
$ws.Range("N11").Value = -3788

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 30000
$ws.Range("I42").Value = 30000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 30000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -29407
$ws.Range("N42").ClearContents()
$ws.Range("H59").Value = 55000
$ws.Range("I59").Value = 60000
$ws.Range("J59").Value = 50000
$ws.Range("K59").Value = 60000
$ws.Range("L59").Value = 50000
$ws.Range("M59").Value = -58855
$ws.Range("N59").Value = -52290
$ws.Range("H62").Value = 4999
$ws.Range("I62").Value = 4998
$ws.Range("K62").Value = 4998
$ws.Range("M62").Value = -4374
$ws.Range("H65").Value = 4999
$ws.Range("I65").Value = 4998
$ws.Range("K65").Value = 24990
$ws.Range("M65").Value = -21870
$ws.Range("H107").Value = 657.5
$ws.Range("I107").Value = 686.55554
$ws.Range("J107").Value = 396
$ws.Range("K107").Value = 686.55554
$ws.Range("L107").Value = 396
$ws.Range("M107").Value = 1233.44446
$ws.Range("N107").Value = -4236
$ws.Range("H132").Value = 8129.4287
$ws.Range("I132").Value = 3726.5
$ws.Range("K132").Value = 11179.5
$ws.Range("M132").Value = -8649.5
$ws.Range("H134").Value = 8957.583000000001
$ws.Range("I134").Value = 4202.2
$ws.Range("K134").Value = 12606.6
$ws.Range("M134").Value = -10071.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 204.44
$ws.Range("I2").Value = 253.95
$ws.Range("J2").Value = 6.4
$ws.Range("K2").Value = 1523.7
$ws.Range("L2").Value = 38.40000000000001
$ws.Range("M2").Value = -1410.7
$ws.Range("N2").Value = -264.4
$ws.Range("H12").Value = 81.5
$ws.Range("I12").Value = 35
$ws.Range("J12").Value = 112.5
$ws.Range("K12").Value = 105
$ws.Range("L12").Value = 337.5
$ws.Range("M12").Value = 68
$ws.Range("N12").Value = -683.5
$ws.Range("H14").Value = 4997.5
$ws.Range("I14").Value = 4997.5
$ws.Range("K14").Value = 14992.5
$ws.Range("M14").Value = -14819.5
$ws.Range("I17").Value = 150
$ws.Range("J17").Value = 1154.2
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 3462.6
$ws.Range("M17").Value = -281
$ws.Range("N17").Value = -3800.6
$ws.Range("H38").Value = 165.5
$ws.Range("I38").Value = 71
$ws.Range("J38").Value = 236.375
$ws.Range("K38").Value = 213
$ws.Range("L38").Value = 709.125
$ws.Range("M38").Value = 134
$ws.Range("N38").Value = -1403.125
$ws.Range("H131").Value = 2170.8572
$ws.Range("I131").Value = 999.5
$ws.Range("J131").Value = 2366.0833
$ws.Range("K131").Value = 2998.5
$ws.Range("L131").Value = 7098.249899999999
$ws.Range("M131").Value = 2041.5
$ws.Range("N131").Value = -17178.2499

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1499
$ws.Range("I80").Value = 1499
$ws.Range("K80").Value = 1499
$ws.Range("M80").Value = -501
$ws.Range("H83").Value = 1499
$ws.Range("I83").Value = 1499
$ws.Range("K83").Value = 7495
$ws.Range("M83").Value = -2503
$ws.Range("H132").Value = 8131.25
$ws.Range("I132").Value = 5113
$ws.Range("J132").Value = 17186
$ws.Range("K132").Value = 15339
$ws.Range("L132").Value = 51558
$ws.Range("M132").Value = -12809
$ws.Range("N132").Value = -56618

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2178.8
$ws.Range("I16").Value = 1498.5
$ws.Range("K16").Value = 1498.5
$ws.Range("M16").Value = -1328.5
$ws.Range("H22").Value = 650
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -505
$ws.Range("H27").Value = 650
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 800
$ws.Range("M27").Value = -693

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 49997.5
$ws.Range("J123").Value = 49997.5
$ws.Range("L123").Value = 49997.5
$ws.Range("N123").Value = -59797.5
$ws.Range("H126").Value = 2213.5557
$ws.Range("I126").Value = 1955.5
$ws.Range("J126").Value = 2729.6667
$ws.Range("K126").Value = 5866.5
$ws.Range("L126").Value = 8189.000100000001
$ws.Range("M126").Value = -3396.5
$ws.Range("N126").Value = -13129.0001
$ws.Range("H132").Value = 10267
$ws.Range("I132").Value = 9200.429
$ws.Range("K132").Value = 27601.287
$ws.Range("M132").Value = -25071.287
$ws.Range("H136").Value = 10250
$ws.Range("I136").Value = 3500
$ws.Range("J136").Value = 17000
$ws.Range("K136").Value = 10500
$ws.Range("L136").Value = 51000
$ws.Range("M136").Value = -7950
$ws.Range("N136").Value = -56100
